$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1 ("电缆" / Cable) - grows from A1:C20 to A1:C25
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert 5 fresh rows before the old row 20 (the 630/240/2 row), so it
# becomes row 25 and we get 5 new blank rows (20-24) to populate.
$ws1.Range("A20:A24").Insert()

# Row 4: 20 / 4 -> 2.5 / 1
$ws1.Cells.Item(4, 2).Value = 2.5

# Row 5: 25 / 6 -> 4 / 1
$ws1.Cells.Item(5, 2).Value = 4

# Row 8: 50 / 16 -> 10 / 1
$ws1.Cells.Item(8, 2).Value = 10

# Row 13: 160 -> 140 / 70 / 1
$ws1.Cells.Item(13, 1).Value = 140

# Row 14: 200 -> 160 / 120 -> 70 / 1
$ws1.Cells.Item(14, 1).Value = 160
$ws1.Cells.Item(14, 2).Value = 70

# Row 15: 250 -> 180 / 150 -> 95 / 1
$ws1.Cells.Item(15, 1).Value = 180
$ws1.Cells.Item(15, 2).Value = 95

# Row 16: 320 -> 200 / 185 -> 95 / 1
$ws1.Cells.Item(16, 1).Value = 200
$ws1.Cells.Item(16, 2).Value = 95

# Row 17: 360 -> 225 / 240 -> 120 / 1
$ws1.Cells.Item(17, 1).Value = 225
$ws1.Cells.Item(17, 2).Value = 120

# Row 18: 400 -> 250 / 120 -> 150 / 2 -> 1
$ws1.Cells.Item(18, 1).Value = 250
$ws1.Cells.Item(18, 2).Value = 150
$ws1.Cells.Item(18, 3).Value = 1

# Row 19: 500 -> 315 / 150 -> 185 / 2 -> 1
$ws1.Cells.Item(19, 1).Value = 315
$ws1.Cells.Item(19, 2).Value = 185
$ws1.Cells.Item(19, 3).Value = 1

# New row 20: 320 / 185 / 1
$ws1.Cells.Item(20, 1).Value = 320
$ws1.Cells.Item(20, 2).Value = 185
$ws1.Cells.Item(20, 3).Value = 1

# New row 21: 350 / 240 / 1
$ws1.Cells.Item(21, 1).Value = 350
$ws1.Cells.Item(21, 2).Value = 240
$ws1.Cells.Item(21, 3).Value = 1

# New row 22: 360 / 240 / 1
$ws1.Cells.Item(22, 1).Value = 360
$ws1.Cells.Item(22, 2).Value = 240
$ws1.Cells.Item(22, 3).Value = 1

# New row 23: 400 / 120 / 2
$ws1.Cells.Item(23, 1).Value = 400
$ws1.Cells.Item(23, 2).Value = 120
$ws1.Cells.Item(23, 3).Value = 2

# New row 24: 500 / 150 / 2
$ws1.Cells.Item(24, 1).Value = 500
$ws1.Cells.Item(24, 2).Value = 150
$ws1.Cells.Item(24, 3).Value = 2

# Row 25 already holds 630 / 240 / 2 (shifted down from old row 20) - unchanged

# Update the view's active cell / selection for sheet 1
$ws1.Range("G13").Select()

# ------------------------------------------------------------------
# Sheet 2 ("电线" / Wire) - grows from A1:C14 to A1:C16
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 4: 20 / 4 -> 2.5 / 1
$ws2.Cells.Item(4, 2).Value = 2.5

# Row 5: 25 / 6 -> 4 / 1
$ws2.Cells.Item(5, 2).Value = 4

# Row 6: 32 / 10 -> 6 / 1
$ws2.Cells.Item(6, 2).Value = 6

# Row 8: 50 / 16 -> 10 / 1
$ws2.Cells.Item(8, 2).Value = 10

# Row 9: 63 / 25 -> 16 / 1
$ws2.Cells.Item(9, 2).Value = 16

# Row 13: 160 -> 140 / 70 / 1
$ws2.Cells.Item(13, 1).Value = 140

# Insert 2 fresh rows before the old row 14 (200/120/1), so it becomes
# row 16 and we get 2 new blank rows (14-15) to populate.
$ws2.Range("A14:A15").Insert()

# New row 14: 160 / 70 / 1
$ws2.Cells.Item(14, 1).Value = 160
$ws2.Cells.Item(14, 2).Value = 70
$ws2.Cells.Item(14, 3).Value = 1

# New row 15: 180 / 95 / 1
$ws2.Cells.Item(15, 1).Value = 180
$ws2.Cells.Item(15, 2).Value = 95
$ws2.Cells.Item(15, 3).Value = 1

# Row 16 (shifted from old row 14): 200 / 120 -> 95 / 1
$ws2.Cells.Item(16, 2).Value = 95

# Update the view's active cell / selection for sheet 2
$ws2.Range("B42").Select()

# Restore sheet 1 as the active/selected tab (it was tabSelected="1"
# originally and only its selection/dimension should have changed).
$ws1.Activate()
